# Fix failing atomic tabOTTR tests due to a breaking change:
# the "data"/"blank" marker rows in the typedFreshBlanks fixture need to be
# realigned - A7 becomes the bare value 1, A8 becomes the "blank" marker
# (previously held by A9), and A9 now reads "data".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 1
$ws.Range("A8").Value = "blank"
$ws.Range("A9").Value = "data"

# Move the active selection to A8, matching the edited cell the author
# was looking at.
$ws.Range("A8").Select() | Out-Null
